$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5:D6").NumberFormat = "@"

$ws.Range("A5").Value = "Saistyphos30"
$ws.Range("C5").Value = "$3.44"
$ws.Range("D5").Value = "$10.33"

$ws.Range("A6").Value = "Iskadarya95"
$ws.Range("C6").Value = "$4.37"
$ws.Range("D6").Value = "$13.10"
